# Auto-generated: Add data for 2023-09-20
# Updates 2023 YTD (column J) and a couple of corrected 2022 (column I) figures
# across "Citywide Totals", "By Neighborhood", and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet = "Citywide Totals"; Cell = "J2"; Value = 5503 },
    @{ Sheet = "Citywide Totals"; Cell = "J3"; Value = 5863 },
    @{ Sheet = "Citywide Totals"; Cell = "I4"; Value = 1304 },
    @{ Sheet = "Citywide Totals"; Cell = "J4"; Value = 1277 },
    @{ Sheet = "Citywide Totals"; Cell = "J5"; Value = 453 },
    @{ Sheet = "Citywide Totals"; Cell = "J6"; Value = 7358 },
    @{ Sheet = "Citywide Totals"; Cell = "I7"; Value = 18481 },
    @{ Sheet = "Citywide Totals"; Cell = "J7"; Value = 20454 },
    @{ Sheet = "By Neighborhood"; Cell = "J2"; Value = 163 },
    @{ Sheet = "By Neighborhood"; Cell = "J7"; Value = 592 },
    @{ Sheet = "By Neighborhood"; Cell = "J8"; Value = 1283 },
    @{ Sheet = "By Neighborhood"; Cell = "J10"; Value = 137 },
    @{ Sheet = "By Neighborhood"; Cell = "J15"; Value = 223 },
    @{ Sheet = "By Neighborhood"; Cell = "J19"; Value = 596 },
    @{ Sheet = "By Neighborhood"; Cell = "J23"; Value = 195 },
    @{ Sheet = "By Neighborhood"; Cell = "J25"; Value = 103 },
    @{ Sheet = "By Neighborhood"; Cell = "J27"; Value = 122 },
    @{ Sheet = "By Neighborhood"; Cell = "J34"; Value = 100 },
    @{ Sheet = "By Neighborhood"; Cell = "J36"; Value = 284 },
    @{ Sheet = "By Neighborhood"; Cell = "J37"; Value = 630 },
    @{ Sheet = "By Neighborhood"; Cell = "J42"; Value = 848 },
    @{ Sheet = "By Neighborhood"; Cell = "J44"; Value = 158 },
    @{ Sheet = "By Neighborhood"; Cell = "J47"; Value = 157 },
    @{ Sheet = "By Neighborhood"; Cell = "J51"; Value = 252 },
    @{ Sheet = "By Neighborhood"; Cell = "J53"; Value = 280 },
    @{ Sheet = "By Neighborhood"; Cell = "J54"; Value = 393 },
    @{ Sheet = "By Neighborhood"; Cell = "J55"; Value = 278 },
    @{ Sheet = "By Neighborhood"; Cell = "J60"; Value = 126 },
    @{ Sheet = "By Neighborhood"; Cell = "J67"; Value = 777 },
    @{ Sheet = "By Neighborhood"; Cell = "J72"; Value = 82 },
    @{ Sheet = "By Neighborhood"; Cell = "J73"; Value = 198 },
    @{ Sheet = "By Neighborhood"; Cell = "J76"; Value = 296 },
    @{ Sheet = "By Neighborhood"; Cell = "J78"; Value = 251 },
    @{ Sheet = "By Neighborhood"; Cell = "J79"; Value = 585 },
    @{ Sheet = "By Neighborhood"; Cell = "J84"; Value = 176 },
    @{ Sheet = "By Neighborhood"; Cell = "J86"; Value = 123 },
    @{ Sheet = "By Neighborhood"; Cell = "J89"; Value = 266 },
    @{ Sheet = "By Neighborhood"; Cell = "J90"; Value = 223 },
    @{ Sheet = "By Neighborhood"; Cell = "J95"; Value = 307 },
    @{ Sheet = "By Neighborhood"; Cell = "J96"; Value = 243 },
    @{ Sheet = "By Neighborhood"; Cell = "I98"; Value = 128 },
    @{ Sheet = "By Neighborhood"; Cell = "J98"; Value = 143 },
    @{ Sheet = "By Neighborhood"; Cell = "J99"; Value = 321 },
    @{ Sheet = "By Neighborhood"; Cell = "I101"; Value = 18481 },
    @{ Sheet = "By Neighborhood"; Cell = "J101"; Value = 20454 },
    @{ Sheet = "West Ridge"; Cell = "J6"; Value = 86 },
    @{ Sheet = "West Ridge"; Cell = "J7"; Value = 243 },
    @{ Sheet = "Auburn Gresham"; Cell = "J2"; Value = 185 },
    @{ Sheet = "Auburn Gresham"; Cell = "J6"; Value = 187 },
    @{ Sheet = "Auburn Gresham"; Cell = "J7"; Value = 592 },
    @{ Sheet = "Uptown"; Cell = "J6"; Value = 78 },
    @{ Sheet = "Uptown"; Cell = "J7"; Value = 266 },
    @{ Sheet = "South Shore"; Cell = "J4"; Value = 56 },
    @{ Sheet = "South Shore"; Cell = "J6"; Value = 251 },
    @{ Sheet = "Logan Square"; Cell = "J6"; Value = 177 },
    @{ Sheet = "Logan Square"; Cell = "J7"; Value = 280 },
    @{ Sheet = "Austin"; Cell = "J2"; Value = 356 },
    @{ Sheet = "Austin"; Cell = "J6"; Value = 432 },
    @{ Sheet = "Austin"; Cell = "J7"; Value = 1283 },
    @{ Sheet = "West Pullman"; Cell = "J3"; Value = 107 },
    @{ Sheet = "West Pullman"; Cell = "J7"; Value = 307 },
    @{ Sheet = "Grand Crossing"; Cell = "J4"; Value = 21 },
    @{ Sheet = "Grand Crossing"; Cell = "J7"; Value = 630 },
    @{ Sheet = "Woodlawn"; Cell = "J3"; Value = 127 },
    @{ Sheet = "Woodlawn"; Cell = "J7"; Value = 321 },
    @{ Sheet = "North Lawndale"; Cell = "J3"; Value = 299 },
    @{ Sheet = "North Lawndale"; Cell = "J5"; Value = 22 },
    @{ Sheet = "North Lawndale"; Cell = "J6"; Value = 203 },
    @{ Sheet = "North Lawndale"; Cell = "J7"; Value = 777 },
    @{ Sheet = "South Deering"; Cell = "J3"; Value = 57 },
    @{ Sheet = "South Deering"; Cell = "J7"; Value = 176 },
    @{ Sheet = "Loop"; Cell = "J6"; Value = 187 },
    @{ Sheet = "Loop"; Cell = "J7"; Value = 393 },
    @{ Sheet = "Chatham"; Cell = "J2"; Value = 150 },
    @{ Sheet = "Chatham"; Cell = "J5"; Value = 23 },
    @{ Sheet = "Chatham"; Cell = "J6"; Value = 217 },
    @{ Sheet = "Chatham"; Cell = "J7"; Value = 596 },
    @{ Sheet = "Irving Park"; Cell = "J6"; Value = 65 },
    @{ Sheet = "Irving Park"; Cell = "J7"; Value = 158 },
    @{ Sheet = "River North"; Cell = "J6"; Value = 164 },
    @{ Sheet = "River North"; Cell = "J7"; Value = 296 },
    @{ Sheet = "Humboldt Park"; Cell = "J2"; Value = 188 },
    @{ Sheet = "Humboldt Park"; Cell = "J3"; Value = 171 },
    @{ Sheet = "Humboldt Park"; Cell = "J4"; Value = 39 },
    @{ Sheet = "Humboldt Park"; Cell = "J6"; Value = 433 },
    @{ Sheet = "Humboldt Park"; Cell = "J7"; Value = 848 },
    @{ Sheet = "Avondale"; Cell = "J3"; Value = 28 },
    @{ Sheet = "Avondale"; Cell = "J7"; Value = 137 },
    @{ Sheet = "Rogers Park"; Cell = "J4"; Value = 27 },
    @{ Sheet = "Rogers Park"; Cell = "J7"; Value = 251 },
    @{ Sheet = "Lower West Side"; Cell = "J4"; Value = 11 },
    @{ Sheet = "Lower West Side"; Cell = "J6"; Value = 137 },
    @{ Sheet = "Lower West Side"; Cell = "J7"; Value = 278 },
    @{ Sheet = "Douglas"; Cell = "J6"; Value = 51 },
    @{ Sheet = "Douglas"; Cell = "J7"; Value = 195 },
    @{ Sheet = "Roseland"; Cell = "J3"; Value = 207 },
    @{ Sheet = "Roseland"; Cell = "J7"; Value = 585 },
    @{ Sheet = "Grand Boulevard"; Cell = "J3"; Value = 91 },
    @{ Sheet = "Grand Boulevard"; Cell = "J6"; Value = 86 },
    @{ Sheet = "Grand Boulevard"; Cell = "J7"; Value = 284 },
    @{ Sheet = "Garfield Ridge"; Cell = "J6"; Value = 37 },
    @{ Sheet = "Garfield Ridge"; Cell = "J7"; Value = 100 },
    @{ Sheet = "East Side"; Cell = "J5"; Value = 4 },
    @{ Sheet = "East Side"; Cell = "J7"; Value = 103 },
    @{ Sheet = "Kenwood"; Cell = "J6"; Value = 73 },
    @{ Sheet = "Kenwood"; Cell = "J7"; Value = 157 },
    @{ Sheet = "Brighton Park"; Cell = "J2"; Value = 65 },
    @{ Sheet = "Brighton Park"; Cell = "J7"; Value = 223 },
    @{ Sheet = "Wicker Park"; Cell = "J2"; Value = 25 },
    @{ Sheet = "Wicker Park"; Cell = "J3"; Value = 22 },
    @{ Sheet = "Wicker Park"; Cell = "I4"; Value = 7 },
    @{ Sheet = "Wicker Park"; Cell = "J6"; Value = 88 },
    @{ Sheet = "Wicker Park"; Cell = "I7"; Value = 128 },
    @{ Sheet = "Wicker Park"; Cell = "J7"; Value = 143 },
    @{ Sheet = "Portage Park"; Cell = "J6"; Value = 65 },
    @{ Sheet = "Portage Park"; Cell = "J7"; Value = 198 },
    @{ Sheet = "Albany Park"; Cell = "J6"; Value = 61 },
    @{ Sheet = "Albany Park"; Cell = "J7"; Value = 163 },
    @{ Sheet = "Edgewater"; Cell = "J6"; Value = 41 },
    @{ Sheet = "Edgewater"; Cell = "J7"; Value = 122 },
    @{ Sheet = "Streeterville"; Cell = "J4"; Value = 66 },
    @{ Sheet = "Streeterville"; Cell = "J7"; Value = 123 },
    @{ Sheet = "Washington Heights"; Cell = "J6"; Value = 64 },
    @{ Sheet = "Washington Heights"; Cell = "J7"; Value = 223 },
    @{ Sheet = "Little Italy, UIC"; Cell = "J6"; Value = 94 },
    @{ Sheet = "Little Italy, UIC"; Cell = "J7"; Value = 252 },
    @{ Sheet = "Morgan Park"; Cell = "J6"; Value = 36 },
    @{ Sheet = "Morgan Park"; Cell = "J7"; Value = 126 },
    @{ Sheet = "Old Town"; Cell = "J6"; Value = 30 },
    @{ Sheet = "Old Town"; Cell = "J7"; Value = 82 }
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    $ws.Range($edit.Cell).Value = $edit.Value
}
